# Add a "language options" menu row (EXIT / BACK / RETOUR / ZURÜCK)
# to the localization table on Sheet1, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "EXIT"
$ws.Range("B5").Value = "BACK"
$ws.Range("C5").Value = "RETOUR"
$ws.Range("D5").Value = "ZURÜCK"

$ws.Range("K19").Select() | Out-Null
